# figureS7_sweep_dynamics/site_transitions.xlsx
# Update the transitional_list sheet with the newly renumbered mutation
# lists for N2.1998B / N2.2002A / N2.2002B (positions shifted after adding
# the new annotated-backbone numbering), and give the first of the three
# rows a slightly shorter, top-aligned wrapped row to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitional_list")

$ws.Range("B3").Value = "V13I,I17L,T19A,V30I,Y40C,N43H,N47T,I62T,E64K,K64N,I77K,K93R,D147N,V149F,I165V,K187R,E199N,K199E, M241V,R249I,I263V,L269I,V313A"
$ws.Range("B4").Value = "L22F,N41D,V50A,M51I,K75R,D86S,K93N,D127G,I149V,S161N,I176M,I194V,D199N,F205L,I257V,V263I,I263V,R264H,K267T,S284F,D309N,H310Y,V312I"
$ws.Range("B5").Value = "I26V,V26I,D41N,N43D,Q49H,L52F,K62T,I73V,L81P,A82T,I257V,I263V,I312T"

# Row 3 now wraps to fewer lines; give it a custom height and top-align
# the wrapped text instead of the default center alignment.
$ws.Rows.Item(3).RowHeight = 47.25
$ws.Range("B3").VerticalAlignment = -4160

$ws.Range("A1:B5").Select()
